$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '24.433.11'
$ws.Range("E2").Value = '  -1.45%  '
Set-TextValue $ws.Range("D3") '1.661.22'
$ws.Range("E3").Value = '  -2.39%  '
Set-TextValue $ws.Range("D4") '1.001'
$ws.Range("E4").Value = '  +0.30%  '
Set-TextValue $ws.Range("D5") '310.98'
$ws.Range("E5").Value = '  -1.15%  '
Set-TextValue $ws.Range("D6") '1.002'
$ws.Range("E6").Value = '  +0.41%  '
Set-TextValue $ws.Range("D7") '0.3908'
$ws.Range("E7").Value = '  -1.19%  '
Set-TextValue $ws.Range("D8") '0.3916'
$ws.Range("E8").Value = '  -2.73%  '
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D9") '50.61'
$ws.Range("E9").Value = '  -4.77%  '
$ws.Range("E10").Value = '  -6.46%  '
$ws.Range("B11").Value = 'BinanceUSD'
$ws.Range("C11").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range("D11") '1.001'
$ws.Range("E11").Value = '  +0.36%  '
Set-TextValue $ws.Range("D12") '0.08569'
$ws.Range("E12").Value = '  -2.98%  '
Set-TextValue $ws.Range("D13") '24.95'
$ws.Range("E13").Value = '  -5.45%  '
Set-TextValue $ws.Range("D14") '7.254'
$ws.Range("E14").Value = '  -3.06%  '
Set-TextValue $ws.Range("D15") '0.00001308'
$ws.Range("E15").Value = '  -3.29%  '
$ws.Range("E16").Value = '  -4.50%  '
Set-TextValue $ws.Range("D17") '1.651.26'
$ws.Range("E17").Value = '  -0.93%  '
Set-TextValue $ws.Range("D18") '94.37'
$ws.Range("E18").Value = '  -1.56%  '
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range("D19") '21.33'
$ws.Range("E19").Value = '  +2.64%  '
$ws.Range("B20").Value = 'TRON'
$ws.Range("C20").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range("D20") '0.06952'
$ws.Range("E20").Value = '  -3.20%  '
Set-TextValue $ws.Range("D21") '7.040'
$ws.Range("E21").Value = '  -4.55%  '
Set-TextValue $ws.Range("D22") '1.003'
$ws.Range("E22").Value = '  +0.64%  '
Set-TextValue $ws.Range("D23") '13.88'
Set-TextValue $ws.Range("D24") '24.441.73'
$ws.Range("E24").Value = '  -1.37%  '
Set-TextValue $ws.Range("D25") '2.378'
$ws.Range("E25").Value = '  +0.75%  '
Set-TextValue $ws.Range("D26") '2.778'
$ws.Range("E26").Value = '  -5.08%  '
Set-TextValue $ws.Range("D27") '22.77'
$ws.Range("E27").Value = '  -2.24%  '
Set-TextValue $ws.Range("D28") '159.62'
$ws.Range("E28").Value = '  -0.77%  '
Set-TextValue $ws.Range("D29") '5.750'
$ws.Range("E29").Value = '  -9.26%  '
Set-TextValue $ws.Range("D30") '145.03'
$ws.Range("E30").Value = '  -0.64%  '
Set-TextValue $ws.Range("D31") '8.148'
$ws.Range("E31").Value = '  -3.47%  '
Set-TextValue $ws.Range("D32") '2.580'
$ws.Range("E32").Value = '  +7.00%  '
Set-TextValue $ws.Range("D33") '1.837.17'
$ws.Range("E33").Value = '  -0.71%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D34") '0.08241'
$ws.Range("E34").Value = '  -4.30%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D35") '1.006'
$ws.Range("E35").Value = '  -2.52%  '
$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D36") '6.875'
$ws.Range("E36").Value = '  -5.20%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D37") '0.02977'
$ws.Range("E37").Value = '  -4.68%  '
Set-TextValue $ws.Range("D38") '0.2777'
$ws.Range("E38").Value = '  -2.64%  '
Set-TextValue $ws.Range("D39") '0.09424'
$ws.Range("E39").Value = '  -0.37%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D40") '1.491'
$ws.Range("E40").Value = '  +1.19%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D41") '10.18'
$ws.Range("E41").Value = '  -5.04%  '
Set-TextValue $ws.Range("D42") '0.7818'
$ws.Range("E42").Value = '  -7.40%  '
Set-TextValue $ws.Range("D43") '13.35'
$ws.Range("E43").Value = '  -5.12%  '
Set-TextValue $ws.Range("D44") '16.27'
$ws.Range("E44").Value = '  -6.55%  '
Set-TextValue $ws.Range("D45") '2.555'
$ws.Range("E45").Value = '  -4.61%  '
Set-TextValue $ws.Range("D46") '0.7041'
$ws.Range("E46").Value = '  -5.56%  '
Set-TextValue $ws.Range("D47") '4.145'
$ws.Range("E47").Value = '  -1.66%  '
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range("D48") '1.002'
$ws.Range("E48").Value = '  +0.51%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D49") '0.08600'
$ws.Range("E49").Value = '  -0.68%  '
Set-TextValue $ws.Range("D50") '1.312'
$ws.Range("E50").Value = '  -4.87%  '
Set-TextValue $ws.Range("D51") '136.05'
$ws.Range("E51").Value = '  -2.43%  '
